$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '30.308.51'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.98%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.921.58'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +0.61%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '0.8119'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +2.52%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '244.20'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.89%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +0.06%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3261'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +3.02%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '27.14'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +2.94%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.07244'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +5.03%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.7928'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +6.54%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.08104'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +1.26%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '1.924.97'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +0.84%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '5.409'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +4.22%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '94.17'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +1.19%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '30.332.51'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +1.09%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '14.30'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +2.43%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '6.083'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +3.66%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '251.07'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +2.14%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.000007856'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +1.45%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '2.183.11'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +1.27%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.15%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '8.057'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +17.86%  '

$ws.Range("E24").Value = '  +0.12%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.1672'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +20.05%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '9.514'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +3.05%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '167.28'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.43%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '19.08'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +0.84%  '

$ws.Range("E29").Value = '  +6.16%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '1.373'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +0.68%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '1.548'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +2.26%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '4.351'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +0.78%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.05682'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +2.69%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '4.139'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +1.44%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.301'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +3.55%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.7454'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +1.52%  '

$ws.Range("E37").Value = '  +0.20%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '2.729'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.31%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.01959'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +1.71%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '2.820'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +1.26%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.4499'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +1.86%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '74.63'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +2.98%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '5.989'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -2.40%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.8565'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +2.31%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.928'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +2.65%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '1.043.24'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +5.50%  '

$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.03%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '103.12'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +2.61%  '

$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '7.637'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +1.10%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '9.908'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +1.39%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '3.087'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +10.33%  '
